# Add two new columns (I0 and IF) to the worksheet, mirroring the
# existing header/data layout already used for columns A-H.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers (row 1) - reuse the exact formatting of the existing header
# cells (e.g. H1) by copying format only, so no new style is created.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows 2-12
$values = @(
    @(5, 6),
    @(6, 7),
    @(7, 9),
    @(3, 7),
    @(4, 6),
    @(2, 6),
    @(7, 9),
    @(8, 8),
    @(6, 7),
    @(5, 6),
    @(6, 6)
)

$row = 2
foreach ($pair in $values) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row = $row + 1
}
